# Rebuild the single-paragraph Bescheid letter: the whole body is one
# <w:r> whose text is broken up by <w:br/> line breaks (char 11 / vertical
# tab in the Range.Text model), so the cleanest, least fragile way to
# reproduce the edit is to rewrite that run's text wholesale rather than
# chase dozens of overlapping short Find/Replace anchors.
$d = $word.ActiveDocument

$segments = @(
    'Einleitung:',
    '',
    'Landratsamt Ortenaukreis',
    'Herrn Franz Konrad',
    'Sachbearbeiter',
    '',
    'Reparatur Ihres Fachwerkhauses in Neuried, Lange Straße 12',
    '',
    'Sehr geehrter Herr Konrad,',
    '',
    'Tenor:',
    '',
    'Es wird angeordnet, dass Sie Ihr Fachwerkhaus mit Biberschwanzdachziegeln reparieren müssen.',
    '',
    'Die Anordnung der sofortigen Vollziehung dieser Maßnahme wird hiermit angeordnet.',
    '',
    'Begründung:',
    '',
    'Sie sind als Forstrat zusammen mit Ihrem Bruder, dem Studenten Georg Konrad, Eigentümer des o.g. Fachwerkhauses. Ihr Haus stammt aus dem Jahre 1865. Das Haus hat die für die damalige Bauweise charakteristischen Wetterdächer sowie die typischen, vorragenden Balkenköpfe an den Erdgeschossbalken. Es ist mit den seinerzeit üblichen Biberschwanz-Dachziegeln gedeckt. Das Fachwerkhaus gehört zu den wenigen voll erhaltenen Exemplaren seiner Art am Oberrhein.',
    '',
    'Durch einen Sturm wurden ca. 50 Biberschwanz-Dachziegel abgedeckt.',
    '',
    'Die Reparaturanordnung stützt sich auf § 1 Abs. 1 in Verbindung mit § 7 Abs. 1 Denkmalschutzgesetz (DSchG). Danach können wir Maßnahmen anordnen, wenn ein Kulturdenkmal gefährdet ist.',
    '',
    'Das Fachwerkhaus ist ein Kulturdenkmal. Es besteht nämlich nach § 2 Abs. 1 DSchG an dem Fachwerkhaus ein öffentliches Erhaltungsinteresse aus heimatgeschichtlichen Gründen. Das Haus stammt aus dem Jahr 1865 und hat die für die damalige Bauweise charakteristischen Wetterdächer sowie die typischen, vorragenden Balkenköpfe an den Erdgeschossbalken und ist mit den seinerzeit üblichen Biberschwanz-Dachziegeln gedeckt.',
    '',
    'Das Kulturdenkmal ist gefährdet, da durch das beeinträchtigte Erscheinungsbild bereits ein Schaden entstanden ist. Der Sturm hat ca. 50 Biberschwanz-Dachziegel abgedeckt.',
    '',
    'Sie sind verpflichtet, die Reparatur durchzuführen.',
    '',
    'Das ergibt sich aus § 7 Polizeigesetz (PolG), da Sie Eigentümer des Fachwerkhauses sind, von dem eine Gefährdung des Denkmals ausgeht.',
    '',
    'Ebenfalls nach denselben Vorschriften verpflichtet, ist Ihr Bruder Georg Konrad, da dieser ebenfalls Eigentümer ist. Als Forstrat sind Sie leistungsfähiger als Ihr Bruder Georg und sind deshalb verpflichtet, die Kosten zu tragen.',
    '',
    'Uns ist hier Ermessen eingeräumt, wobei wir als gesetzliche Grenze im Sinne von § 40 Landesverwaltungsverfahrensgesetz die Verhältnismäßigkeit nach Art. 20 Abs. 3 GG beachtet haben.',
    '',
    'Die Anordnung das Dach mit Biberschwanz-Dachziegeln zu reparieren ist angemessen, um die Ansehnlichkeit des Denkmals herbeizuführen. Diese Anordnung ist auch erforderlich, da es kein milderes gleichgeeignetes Mittel gibt. Ethanitplatten wären zwar kostengünstiger, könnten aber nicht die Ansehnlichkeit wieder herstellen. Schließlich ist die Reparaturanordnung auch angemessen.',
    '',
    'Ihr Nachteil steht nicht außer Verhältnis zu den Vorteilen für die Allgemeinheit. Das Interesse am Erhalt des Denkmals gewichtet durch Art. 3 C Landes Verfassung überwiegt Ihr Interesse bezüglich finanzieller Belastung und der Eigentumsrechte gewichtet durch Art. 14 GG.',
    '',
    'Die Durchführung der Maßnahme ist Ihnen auch ohne die Mitwirkung Ihres Bruders Georg möglich ist. An sich hätte er als Miterbe nach § 2038 BGB mitzubestimmen, weil es sich bei der Reparaturanordnung jedoch um eine Notmaßnahme handelt, ist diese Mitbestimmung nicht notwendig. Im Hinblick auf den Verkaufswert des Fachwerkhauses ist es als Notmaßnahme wirtschaftlich vernünftig das Dach mit Biberschwanz-Dachziegeln zu reparieren.',
    '',
    'Als Rechtsgrundlage für die Anordnung dient § 80 Abs. 2 Satz 1 Nr. 4. Landesverwaltungsverfahrensgesetz.',
    '',
    'Rechtsbehelfsbelehrung:',
    '',
    'Gegen die Dachdeckungsanordnung können Sie innerhalb eines Monats nach Bekanntgabe bei dem Landratsamt Ortenaukreis, Badstraße 22, 77652 Offenburg Widerspruch einlegen. (§37 (6) LVwVfG (§70 VwGO))',
    '',
    'Gegen die Anordnung der sofortigen Vollziehung können Sie den Antrag auf Wiederherstellung der aufschiebenden Wirkung beim Verwaltungsgericht Freiburg, Herbstburgerstraße 115, 79104 Freiburg stellen. (§ 80 Abs. 5 VwGO)',
    '',
    'Unterschrift mit Grußformel:',
    '',
    'Mit freundlichen Grüßen,',
    '',
    'Lisa Brunzel'
)
$newText = $segments -join [char]11

$d.Content.Text = $newText
